# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style (bold, border, centered) from the existing header cell G1
# onto the new header cell H1 before filling in its value.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Header text for the new column.
$ws.Range("H1").Value = "Save"

# New "Save" column values for each data row.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 0
$ws.Range("H8").Value = 1
